# Applies the 02-01-2024 scraper-refresh edit:
#  1) Re-orders the match rows within a handful of same-matchday groups
#     (betexplorer re-scrape changed the listing order of matches that
#     share a kickoff date/time).
#  2) Appends 9 newly-scraped matches (rows 83-91) for 21/12, 22/12 and
#     28/12/2023 kickoffs.
#  3) Grows the sheet's used range from A1:V82 to A1:V91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Part 1: re-order rows within same-kickoff-time groups.
# Columns A-E (Indice/pais/torneio/temporada/data_partida) are identical
# for every row in a group, so only F:V (home..url) need to move.
# Mapping is {NewRowNumber -> SourceRowNumber} read from the *original*
# (pre-edit) sheet.
# ---------------------------------------------------------------------
$rowMap = @(
    @{New=42; Old=44}, @{New=43; Old=45}, @{New=44; Old=43}, @{New=45; Old=42},
    @{New=56; Old=60}, @{New=57; Old=59}, @{New=58; Old=57}, @{New=59; Old=56}, @{New=60; Old=58},
    @{New=61; Old=64}, @{New=62; Old=65}, @{New=64; Old=62}, @{New=65; Old=61},
    @{New=67; Old=68}, @{New=68; Old=69}, @{New=69; Old=67},
    @{New=73; Old=74}, @{New=74; Old=75}, @{New=75; Old=76}, @{New=76; Old=77}, @{New=77; Old=73},
    @{New=79; Old=80}, @{New=80; Old=81}, @{New=81; Old=79}
)

$firstCol = 6   # F
$lastCol  = 22  # V

# Snapshot every source row's F:V values BEFORE any writes happen, so
# writes never clobber a value that is still needed as a source later.
$snapshot = @{}
foreach ($entry in $rowMap) {
    $oldRow = $entry.Old
    if (-not $snapshot.ContainsKey($oldRow)) {
        $vals = @()
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $vals += ,($ws.Cells.Item($oldRow, $c).Value2)
        }
        $snapshot[$oldRow] = $vals
    }
}

foreach ($entry in $rowMap) {
    $newRow = $entry.New
    $vals = $snapshot[$entry.Old]
    $i = 0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($newRow, $c).Value2 = $vals[$i]
        $i++
    }
}

# ---------------------------------------------------------------------
# Part 2: append the 9 newly scraped matches as rows 83-91.
# ---------------------------------------------------------------------
$newRows = @(
    @{A=82; E=45283.66666666666; F='Annan';          G=0; H='Falkirk';       I=3; J=7.07;  K='21/12/2023 09:12'; L=10.78; M='23/12/2023 15:40'; N=4.77;  O='21/12/2023 09:12'; P=6.56;  Q='23/12/2023 15:40'; R=1.35;  S='21/12/2023 09:12'; T=1.23;  U='23/12/2023 15:40'; V='https://www.betexplorer.com/football/scotland/league-one/annan-falkirk/fBP3mkTk/'},
    @{A=83; E=45283.66666666666; F='Alloa';           G=2; H='Montrose';      I=2; J=2.44;  K='21/12/2023 09:12'; L=2.02;  M='23/12/2023 15:50'; N=3.24;  O='21/12/2023 09:12'; P=3.35;  Q='23/12/2023 15:57'; R=2.67;  S='21/12/2023 09:12'; T=3.77;  U='23/12/2023 15:57'; V='https://www.betexplorer.com/football/scotland/league-one/alloa-montrose/t8LalVDq/'},
    @{A=84; E=45283.66666666666; F='Hamilton';        G=1; H='Kelty Hearts';  I=1; J=1.29;  K='21/12/2023 09:12'; L=1.39;  M='23/12/2023 15:53'; N=5.1;   O='21/12/2023 09:12'; P=4.78;  Q='23/12/2023 15:53'; R=8.52;  S='21/12/2023 09:12'; T=7.77;  U='23/12/2023 15:53'; V='https://www.betexplorer.com/football/scotland/league-one/hamilton-kelty-hearts/SONBoTc2/'},
    @{A=85; E=45283.66666666666; F='Stirling';        G=1; H='Queen of South';I=1; J=2.34;  K='21/12/2023 09:12'; L=2.67;  M='23/12/2023 15:56'; N=3.26;  O='21/12/2023 09:12'; P=3.44;  Q='23/12/2023 15:55'; R=2.78;  S='21/12/2023 09:12'; T=2.54;  U='23/12/2023 15:55'; V='https://www.betexplorer.com/football/scotland/league-one/stirling-queen-of-south/69ghfL12/'},
    @{A=86; E=45283.6875;        F='Cove Rangers';    G=7; H='Edinburgh City';I=2; J=1.24;  K='22/12/2023 04:13'; L=1.11;  M='23/12/2023 15:08'; N=5.53;  O='22/12/2023 04:13'; P=9.6;   Q='23/12/2023 15:34'; R=9.65;  S='22/12/2023 04:13'; T=20.26; U='23/12/2023 15:34'; V='https://www.betexplorer.com/football/scotland/league-one/cove-rangers-edinburgh-city/lYO7n9rd/'},
    @{A=87; E=45290.66666666666; F='Queen of South';  G=2; H='Annan';         I=1; J=1.69;  K='28/12/2023 09:12'; L=1.91;  M='30/12/2023 15:17'; N=3.82;  O='28/12/2023 09:12'; P=3.77;  Q='30/12/2023 15:17'; R=4.17;  S='28/12/2023 09:12'; T=3.72;  U='30/12/2023 15:17'; V='https://www.betexplorer.com/football/scotland/league-one/queen-of-south-annan/zVWznc9r/'},
    @{A=88; E=45290.66666666666; F='Falkirk';         G=5; H='Stirling';      I=0; J=1.22;  K='28/12/2023 09:12'; L=1.14;  M='30/12/2023 15:28'; N=5.75;  O='28/12/2023 09:12'; P=8.06;  Q='30/12/2023 15:28'; R=10.7;  S='28/12/2023 09:12'; T=17.91; U='30/12/2023 15:28'; V='https://www.betexplorer.com/football/scotland/league-one/falkirk-stirling/vNi0haWE/'},
    @{A=89; E=45290.66666666666; F='Edinburgh City';  G=2; H='Hamilton';      I=5; J=19.53; K='29/12/2023 14:42'; L=35.63; M='30/12/2023 15:57'; N=9.6;   O='29/12/2023 14:42'; P=16.4;  Q='30/12/2023 15:57'; R=1.08;  S='29/12/2023 14:42'; T=1.04;  U='30/12/2023 15:55'; V='https://www.betexplorer.com/football/scotland/league-one/edinburgh-city-hamilton/p0hdguH8/'},
    @{A=90; E=45290.66666666666; F='Montrose';        G=1; H='Cove Rangers';  I=1; J=3.39;  K='28/12/2023 09:12'; L=3.43;  M='30/12/2023 15:59'; N=3.56;  O='28/12/2023 09:12'; P=3.88;  Q='30/12/2023 15:59'; R=1.93;  S='28/12/2023 09:12'; T=1.97;  U='30/12/2023 15:43'; V='https://www.betexplorer.com/football/scotland/league-one/montrose-cove-rangers/YLp9jw1R/'}
)

$r = 83
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2  = $row.A           # Indice
    $ws.Cells.Item($r, 2).Value2  = 'scotland'        # pais
    $ws.Cells.Item($r, 3).Value2  = 'league-one'       # torneio
    $ws.Cells.Item($r, 4).Value2  = '2023-2024'        # temporada
    $ws.Cells.Item($r, 5).Value2  = $row.E             # data_partida
    $ws.Cells.Item($r, 6).Value2  = $row.F             # home
    $ws.Cells.Item($r, 7).Value2  = $row.G             # home_ft_gols
    $ws.Cells.Item($r, 8).Value2  = $row.H             # away
    $ws.Cells.Item($r, 9).Value2  = $row.I             # away_ft_gols
    $ws.Cells.Item($r, 10).Value2 = $row.J             # home_opening_odds
    $ws.Cells.Item($r, 11).Value2 = $row.K             # home_opening_data_hora
    $ws.Cells.Item($r, 12).Value2 = $row.L             # home_closing_odds
    $ws.Cells.Item($r, 13).Value2 = $row.M             # home_closing_data_hora
    $ws.Cells.Item($r, 14).Value2 = $row.N             # draw_opening_odds
    $ws.Cells.Item($r, 15).Value2 = $row.O             # draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value2 = $row.P             # draw_closing_odds
    $ws.Cells.Item($r, 17).Value2 = $row.Q             # draw_closing_data_hora
    $ws.Cells.Item($r, 18).Value2 = $row.R             # away_opening_odds
    $ws.Cells.Item($r, 19).Value2 = $row.S             # away_opening_data_hora
    $ws.Cells.Item($r, 20).Value2 = $row.T             # away_closing_odds
    $ws.Cells.Item($r, 21).Value2 = $row.U             # away_closing_data_hora
    $ws.Cells.Item($r, 22).Value2 = $row.V             # url_partida

    # Match the formatting used by every other data row: bordered/centred
    # bold style on Indice (A) and the custom datetime format on
    # data_partida (E). (.Style assignment is a no-op on this host, so
    # copy/paste-special the formats from an existing row instead.)
    $ws.Cells.Item(2, 1).Copy() | Out-Null
    $ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $ws.Cells.Item(2, 5).Copy() | Out-Null
    $ws.Cells.Item($r, 5).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

    $r++
}
